$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.282.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'219.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.82%  "

$ws.Range("D6").Value = "'0.5292"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").Value = "'1.008"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("E8").Value = "  +0.79%  "

$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("E10").Value = "  +2.49%  "

$ws.Range("D11").Value = "'0.07832"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").Value = "'4.524"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").Value = "'1.679.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.12%  "

$ws.Range("D14").Value = "'1.896.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("D15").Value = "'0.5594"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "'0.0₅8099"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("D17").Value = "'65.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").Value = "'26.297.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").Value = "'1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("D20").Value = "'4.721"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.75%  "

$ws.Range("D21").Value = "'200.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.32%  "

$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("D23").Value = "'6.060"
$ws.Range("D23").Style = "Normal"

$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("D25").Value = "'146.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.60%  "

$ws.Range("D26").Value = "'0.1219"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").Value = "'7.227"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("D29").Value = "'1.529"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.35%  "

$ws.Range("D30").Value = "'0.05882"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.43%  "

$ws.Range("D31").Value = "'1.283"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.42%  "

$ws.Range("D32").Value = "'3.506"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.97%  "

$ws.Range("D33").Value = "'3.333"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.40%  "

$ws.Range("D34").Value = "'1.597"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.91%  "

$ws.Range("D35").Value = "'0.9637"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("D36").Value = "'2.821"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("D37").Value = "'2.428"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("E39").Value = "  +0.46%  "

$ws.Range("D40").Value = "'5.979"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.43%  "

$ws.Range("D41").Value = "'1.079.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.62%  "

$ws.Range("D42").Value = "'0.8587"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.92%  "

$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("D44").Value = "'102.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.94%  "

$ws.Range("D45").Value = "'1.808.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("D46").Value = "'58.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.70%  "

$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").Value = "'0.4415"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.97%  "

$ws.Range("D49").Value = "'8.043"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.69%  "

$ws.Range("E50").Value = "  -5.74%  "

$ws.Range("D51").Value = "'0.05143"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
